# Horarios actualizados Linea 141 - scrape refresh at 20:12:18
# Updates Hora_Scrap/Hora_Llegada/Linea/Minutos rows that moved position
# (the feed is kept sorted by estimated arrival time, so a refreshed
# "minutes remaining" estimate can re-order same-minute rows) and
# appends newly scraped rows at the bottom of each sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912" ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 20:12:18"
$ws.Range("A3").Value = "Total filas: 350"
$ws.Range("C49").Value = "11_ETCHEVERRY"
$ws.Range("C50").Value = "15_ABASTO"
$ws.Range("C51").Value = "11_ETCHEVERRY"
$ws.Range("C52").Value = "15_ABASTO"
$ws.Range("C64").Value = "23_HERNANDEZ"
$ws.Range("C65").Value = "215B_EL PATO"
$ws.Range("A75").Value = "08:40:59"
$ws.Range("C75").Value = "15X38_ABASTO"
$ws.Range("D75").Value = 37
$ws.Range("A76").Value = "08:30:14"
$ws.Range("C76").Value = "27_EL RETIRO"
$ws.Range("D76").Value = 47
$ws.Range("A77").Value = "08:52:33"
$ws.Range("C77").Value = "14_ABASTO"
$ws.Range("D77").Value = 25
$ws.Range("A112").Value = "10:56:01"
$ws.Range("C112").Value = "81_EL PELIGRO"
$ws.Range("D112").Value = 5
$ws.Range("A113").Value = "09:23:52"
$ws.Range("C113").Value = "10_OLMOS"
$ws.Range("D113").Value = 98
$ws.Range("C140").Value = "15_ABASTO"
$ws.Range("C141").Value = "16_P MOR-SANTA ANA"
$ws.Range("C149").Value = "23_HERNANDEZ"
$ws.Range("C150").Value = "14_ABASTO"
$ws.Range("A153").Value = "10:56:01"
$ws.Range("C153").Value = "27_EL RETIRO"
$ws.Range("D153").Value = 98
$ws.Range("A154").Value = "12:33:54"
$ws.Range("C154").Value = "15_ABASTO"
$ws.Range("D154").Value = 1
$ws.Range("C155").Value = "27_EL RETIRO"
$ws.Range("C156").Value = "23_HERNANDEZ"
$ws.Range("C157").Value = "27_EL RETIRO"
$ws.Range("C158").Value = "23_HERNANDEZ"
$ws.Range("A174").Value = "13:14:41"
$ws.Range("C174").Value = "14_ABASTO"
$ws.Range("D174").Value = 18
$ws.Range("A175").Value = "12:33:54"
$ws.Range("C175").Value = "215A_EL PATO"
$ws.Range("D175").Value = 59
$ws.Range("A224").Value = "14:53:58"
$ws.Range("C224").Value = "16_P MOR-SANTA ANA"
$ws.Range("D224").Value = 59
$ws.Range("A225").Value = "14:33:43"
$ws.Range("C225").Value = "27_EL RETIRO"
$ws.Range("D225").Value = 79
$ws.Range("A227").Value = "15:47:47"
$ws.Range("C227").Value = "16_P MOR-SANTA ANA"
$ws.Range("D227").Value = 6
$ws.Range("A229").Value = "14:47:05"
$ws.Range("C229").Value = "27_EL RETIRO"
$ws.Range("D229").Value = 66
$ws.Range("C230").Value = "15X38_ABASTO"
$ws.Range("A257").Value = "16:44:07"
$ws.Range("C257").Value = "23_HERNANDEZ"
$ws.Range("D257").Value = 23
$ws.Range("A258").Value = "16:52:32"
$ws.Range("C258").Value = "16_P MOR-SANTA ANA"
$ws.Range("D258").Value = 15
$ws.Range("C283").Value = "15_ABASTO"
$ws.Range("C284").Value = "16_P MOR-SANTA ANA"
$ws.Range("A292").Value = "17:51:15"
$ws.Range("C292").Value = "14_ABASTO"
$ws.Range("D292").Value = 39
$ws.Range("A293").Value = "17:38:13"
$ws.Range("C293").Value = "23_HERNANDEZ"
$ws.Range("D293").Value = 52
$ws.Range("C295").Value = "15X38_ABASTO"
$ws.Range("C296").Value = "23_HERNANDEZ"
$ws.Range("A309").Value = "18:52:36"
$ws.Range("C309").Value = "27_EL RETIRO"
$ws.Range("D309").Value = 24
$ws.Range("A310").Value = "19:12:11"
$ws.Range("C310").Value = "17_ROMERO"
$ws.Range("D310").Value = 4
$ws.Range("A331").Value = "20:12:18"
$ws.Range("D331").Value = 0
$ws.Range("A335").Value = "20:12:18"
$ws.Range("D335").Value = 19
$ws.Range("A337").Value = "20:12:18"
$ws.Range("D337").Value = 22
$ws.Range("A339").Value = "20:12:18"
$ws.Range("D339").Value = 34
$ws.Range("A340").Value = "20:12:18"
$ws.Range("D340").Value = 35
$ws.Range("C341").Value = "17_ROMERO"
$ws.Range("C342").Value = "215B_EL PATO"
$ws.Range("A344").Value = "20:12:18"
$ws.Range("D344").Value = 43
$ws.Range("A345").Value = "20:12:18"
$ws.Range("D345").Value = 45
$ws.Range("A347").Value = "20:12:18"
$ws.Range("D347").Value = 55
$ws.Range("A348").Value = "20:12:18"
$ws.Range("D348").Value = 58
$ws.Range("A349").Value = "20:12:18"
$ws.Range("D349").Value = 76
$ws.Range("A351").Value = "20:12:18"
$ws.Range("B351").Value = "21:31"
$ws.Range("C351").Value = "16_SANTA ANA"
$ws.Range("D351").Value = 79
$ws.Range("A352").Value = "20:12:18"
$ws.Range("B352").Value = "21:34"
$ws.Range("C352").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D352").Value = 82
$ws.Range("A353").Value = "20:12:18"
$ws.Range("B353").Value = "21:45"
$ws.Range("C353").Value = "23_HERNANDEZ"
$ws.Range("D353").Value = 93
$ws.Range("E353").Value = "LP1912"
$ws.Range("A354").Value = "20:12:18"
$ws.Range("B354").Value = "21:46"
$ws.Range("C354").Value = "14X44_ABASTO"
$ws.Range("D354").Value = 94
$ws.Range("E354").Value = "LP1912"
$ws.Range("A355").Value = "20:12:18"
$ws.Range("B355").Value = "22:04"
$ws.Range("C355").Value = "15_ABASTO"
$ws.Range("D355").Value = 112
$ws.Range("E355").Value = "LP1912"

# --- Sheet "LP1912-215" ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 20:12:18"
$ws.Range("A52").Value = "20:12:18"
$ws.Range("D52").Value = 35

# --- Sheet "6203-6173" ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 20:12:18"
$ws.Range("A3").Value = "Total filas: 45"
$ws.Range("A48").Value = "20:12:18"
$ws.Range("D48").Value = 40
$ws.Range("A50").Value = "20:12:18"
$ws.Range("B50").Value = "21:30"
$ws.Range("C50").Value = "215C_LA PLATA"
$ws.Range("D50").Value = 78
$ws.Range("E50").Value = "L6203"
